# Automatische test-sync: 2025-08-05 19:44:50
# Adds Testmail #11 ("Weten jullie al iets over mijn retour?") as a new
# row at the bottom of the "Logs" sheet, extends the conditional
# formatting ranges to cover the new row, and bumps the "Retour /
# Terugbetaling" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$newRow = 52

$logs.Cells.Item($newRow, 1).Value = "Weten jullie al iets over mijn retour?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #11: Weten jullie al iets over mijn retour?"
$logs.Cells.Item($newRow, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nBedankt voor je e-mail. Om je vraag over je retour te kunnen beantwoorden, hebben we wat meer informatie nodig. Zou je ons alsjeblieft je ordernummer en/of trackingnummer kunnen geven? Hiermee kunnen we het proces van je retourzending bekijken en je van de juiste informatie voorzien.`nAlvast bedankt voor je medewerking.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 19:44:30"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# The reply text in column E contains embedded line breaks; undo the
# row-height bump that a multi-line value would otherwise trigger so the
# row keeps using the sheet's default (non-custom) height, same as every
# other row in the log.
$logs.Rows.Item($newRow).AutoFit()

# Extend each conditional-formatting block (one per column) so it keeps
# covering rows 2..52 instead of just 2..51.
$cfColumns = "D", "G", "H", "I", "J"
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "51")
    $newRange = $logs.Range($col + "2:" + $col + "52")
    $conditions = $oldRange.FormatConditions
    for ($i = 1; $i -le $conditions.Count; $i++) {
        $conditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Dashboard: one more "Retour / Terugbetaling" mail was logged.
$dash.Range("B6").Value = 4
